# Adding Tik Tok domains to the "Exact List" blocklist sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exact List")

# New domains to block (column B), each paired with 0.0.0.0 in column A,
# appended after the existing last row (1389).
$domains = @(
    'v16a.tiktokcdn.com',
    'p16-tiktokcdn-com.akamaized.net',
    'log.tiktokv.com',
    'ib.tiktokv.com',
    'api-h2.tiktokv.com',
    'v16m.tiktokcdn.com',
    'api.tiktokv.com',
    'v19.tiktokcdn.com',
    'mon.musical.ly',
    'api2-16-h2.musical.ly',
    'api2.musical.ly',
    'log2.musical.ly',
    'api2-21-h2.musical.ly',
    'api-h2.tiktokv.com',
    'api.tiktokv.com',
    'api21-h2.tiktokv.com',
    'ib.tiktokv.com',
    'm.tiktok.com',
    'muscdn.com',
    'tiktokcdn.com',
    'tiktokcdn.com.c.worldfcdn.com',
    'tiktok.com',
    '*.tiktok.com',
    'www.tiktok.com',
    'api-h2.tiktokv.com',
    'api21-h2.tiktokv.com',
    'ns-440.awsdns-55.com',
    'ns-722.awsdns-26.net',
    'ns-1475.awsdns-56.org',
    'ns-1574.awsdns-04.co.uk',
    'tiktokcdn-com.akamaized.net',
    'api-h2.tiktokv.com',
    'api21-h2.tiktokv.com',
    'ns-440.awsdns-55.com',
    'ns-722.awsdns-26.net',
    'ns-1475.awsdns-56.org',
    'ns-1574.awsdns-04.co.uk',
    'tiktokcdn-com.akamaized.net',
    'm.tiktok.com',
    't.tiktok.com',
    'tiktok.com',
    'tiktok.org',
    'tiktokcdn.com',
    'tiktokv.com',
    'vt.tiktok.com',
    'v19.tiktokcdn.com',
    'v16.tiktokcdn.com',
    'api.tiktokv.com',
    'mon.tiktokv.com',
    'ib.tiktokv.com',
    'p16-tiktokcdn-com.akamaized.net',
    'v16-tiktokcdn-com.akamaized.net',
    'musical.ly',
    'muscdn.com',
    'akamai.net',
    'isnssdk.com',
    'www.tiktok.com'
)

$startRow = 1390
for ($i = 0; $i -lt $domains.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "0.0.0.0"
    $ws.Cells.Item($r, 2).Value = $domains[$i]
}

$endRow = $startRow + $domains.Length - 1

# Update the sheet view to reflect the new selection / scroll position,
# matching what Excel would record after scrolling to and selecting the
# newly appended rows.
[void]$ws.Select()
$excel.ActiveWindow.ScrollRow = 1361
[void]$ws.Range("A1389:A" + $endRow).Select()
